$wb = $excel.ActiveWorkbook

# --- Rename sheets (workbook.xml <sheet name="..."/>) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16502911377317884"
$wb.Worksheets.Item(2).Name = "NB_TO-1650291139894714"
$wb.Worksheets.Item(3).Name = "RS_TO-16502911398967123"
$wb.Worksheets.Item(4).Name = "TOL_TO-16502911399687996"
$wb.Worksheets.Item(5).Name = "vSAT_TO-1650291140039658"

# --- Sheet1 (GNG_TO) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16502911376807868.csv"
$ws1.Range("B3").Value = "GNG_stims-16502911376997888.csv"
$ws1.Range("B4").Value = "go_stims-1650291137700786.csv"
$ws1.Range("B5").Value = "GNG_stims-16502911377307863.csv"

# --- Sheet2 (NB_TO) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-1650291138372812.csv"
$ws2.Range("B3").Value = "TB-16502911387701876.csv"
$ws2.Range("B4").Value = "OB-16502911387418919.csv"
$ws2.Range("B5").Value = "TB-16502911393309903.csv"
$ws2.Range("B6").Value = "TB-1650291139870063.csv"
$ws2.Range("B7").Value = "ZB-match_6-16502911378957844.csv"
$ws2.Range("B8").Value = "OB-16502911382047856.csv"
$ws2.Range("B9").Value = "ZB-match_6-16502911377437847.csv"
$ws2.Range("B10").Value = "ZB-match_1-16502911379927857.csv"

# --- Sheet3 (RS_TO) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet4 (TOL_TO) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16502911399205165.csv"
$ws4.Range("B3").Value = "ZM_stims-16502911398977923.csv"
$ws4.Range("B4").Value = "MM_stims-1650291139951909.csv"
$ws4.Range("B5").Value = "ZM_stims-16502911399215195.csv"
$ws4.Range("B6").Value = "MM_stims-16502911399680593.csv"
$ws4.Range("B7").Value = "ZM_stims-16502911399529433.csv"

# --- Sheet5 (vSAT_TO) ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-1650291139971799.csv"
$ws5.Range("B3").Value = "vSAT_stims-16502911399998226.csv"
$ws5.Range("B4").Value = "SAT_stims-16502911399838853.csv"
$ws5.Range("B5").Value = "vSAT_stims-16502911400245311.csv"
